$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Ultima Actualizacion" timestamp column (D) down one block and
# stamp the newest refresh time on top, mirroring the automatic updater run
# "Actualizar 03-08-2021 12-52-33".

$ws.Range("D30:D43").Value = 44263.49352696759
$ws.Range("D16:D29").Value = 44263.51489300926
$ws.Range("D2:D15").Value = 44263.53626992822
